$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3300
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3052
$ws.Range("H67").Value = 3300
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -2442
$ws.Range("H82").Value = 3540
$ws.Range("I82").Value = 2726.6667
$ws.Range("J82").Value = 5980
$ws.Range("K82").Value = 8180.000100000001
$ws.Range("L82").Value = 17940
$ws.Range("M82").Value = -7774.000100000001
$ws.Range("N82").Value = -18752
$ws.Range("H85").Value = 3540
$ws.Range("I85").Value = 2726.6667
$ws.Range("J85").Value = 5980
$ws.Range("K85").Value = 8180.000100000001
$ws.Range("L85").Value = 17940
$ws.Range("M85").Value = -6776.000100000001
$ws.Range("N85").Value = -20748
$ws.Range("H132").Value = 3016.9714
$ws.Range("I132").Value = 2862.3438
$ws.Range("J132").Value = 4666.3335
$ws.Range("K132").Value = 8587.0314
$ws.Range("L132").Value = 13999.0005
$ws.Range("M132").Value = -6057.0314
$ws.Range("N132").Value = -19059.0005
$ws.Range("H138").Value = 2399.8333
$ws.Range("J138").Value = 3047.2593
$ws.Range("L138").Value = 9141.777900000001
$ws.Range("N138").Value = -19421.7779

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 0
$ws.Range("L4").ClearContents()
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = -282
$ws.Range("H23").Value = 30001.8
$ws.Range("I23").Value = 45006
$ws.Range("J23").Value = 19999
$ws.Range("K23").Value = 45006
$ws.Range("L23").Value = 19999
$ws.Range("M23").Value = -44747
$ws.Range("N23").Value = -20517
$ws.Range("H44").Value = 21803.3
$ws.Range("J44").Value = 21803.3
$ws.Range("L44").Value = 21803.3
$ws.Range("N44").Value = -22779.3
$ws.Range("H55").Value = 27999.334
$ws.Range("J55").Value = 27999.334
$ws.Range("L55").Value = 27999.334
$ws.Range("N55").Value = -28629.334
$ws.Range("H61").Value = 13891099
$ws.Range("I61").Value = 20834680
$ws.Range("J61").Value = 3939.25
$ws.Range("K61").Value = 20834680
$ws.Range("L61").Value = 3939.25
$ws.Range("M61").Value = -20834468
$ws.Range("N61").Value = -4363.25
$ws.Range("H80").Value = 21999
$ws.Range("J80").Value = 21999
$ws.Range("L80").Value = 21999
$ws.Range("N80").Value = -23995
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H83").Value = 21999
$ws.Range("J83").Value = 21999
$ws.Range("L83").Value = 65997
$ws.Range("N83").Value = -75981
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H136").Value = 13891099
$ws.Range("I136").Value = 20834680
$ws.Range("J136").Value = 3939.25
$ws.Range("K136").Value = 62504040
$ws.Range("L136").Value = 11817.75
$ws.Range("M136").Value = -62501490
$ws.Range("N136").Value = -16917.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19999
$ws.Range("J35").Value = 19999
$ws.Range("L35").Value = 19999
$ws.Range("N35").Value = -20619
$ws.Range("H134").Value = 3845.5454
$ws.Range("I134").Value = 4183.8184
$ws.Range("J134").Value = 3507.2727
$ws.Range("K134").Value = 12551.4552
$ws.Range("L134").Value = 10521.8181
$ws.Range("M134").Value = -10016.4552
$ws.Range("N134").Value = -15591.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4778.4653
$ws.Range("I31").Value = 1113.0454
$ws.Range("J31").Value = 8618.429
$ws.Range("K31").Value = 1113.0454
$ws.Range("L31").Value = 8618.429
$ws.Range("M31").Value = -818.0454
$ws.Range("N31").Value = -9208.429
$ws.Range("H34").Value = 4778.4653
$ws.Range("I34").Value = 1113.0454
$ws.Range("J34").Value = 8618.429
$ws.Range("K34").Value = 1113.0454
$ws.Range("L34").Value = 8618.429
$ws.Range("M34").Value = -911.0454
$ws.Range("N34").Value = -9022.429
$ws.Range("H68").Value = 28166
$ws.Range("J68").Value = 28166
$ws.Range("L68").Value = 28166
$ws.Range("N68").Value = -29664
$ws.Range("H71").Value = 28166
$ws.Range("J71").Value = 28166
$ws.Range("L71").Value = 84498
$ws.Range("N71").Value = -91986
$ws.Range("H87").Value = 65000
$ws.Range("J87").Value = 65000
$ws.Range("L87").Value = 65000
$ws.Range("N87").Value = -67372
$ws.Range("H90").Value = 65000
$ws.Range("J90").Value = 65000
$ws.Range("L90").Value = 195000
$ws.Range("N90").Value = -206856
$ws.Range("H141").Value = 172581.5
$ws.Range("J141").Value = 168664.58
$ws.Range("L141").Value = 168664.58
$ws.Range("N141").Value = -179024.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1198.8552
$ws.Range("J68").Value = 1398.4807
$ws.Range("L68").Value = 4195.4421
$ws.Range("N68").Value = -5817.4421
$ws.Range("H71").Value = 1198.8552
$ws.Range("J71").Value = 1398.4807
$ws.Range("L71").Value = 12586.3263
$ws.Range("N71").Value = -20698.3263
$ws.Range("H92").Value = 800
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H97").Value = 598.5
$ws.Range("I97").Value = 595
$ws.Range("J97").Value = 602
$ws.Range("K97").Value = 1785
$ws.Range("L97").Value = 1806
$ws.Range("M97").Value = -1289
$ws.Range("N97").Value = -2798
$ws.Range("H113").Value = 833.64105
$ws.Range("I113").Value = 529
$ws.Range("K113").Value = 1587
$ws.Range("M113").Value = 583
$ws.Range("H131").Value = 1142.7097
$ws.Range("I131").Value = 1004
$ws.Range("J131").Value = 1183.1666
$ws.Range("K131").Value = 3012
$ws.Range("L131").Value = 3549.4998
$ws.Range("M131").Value = 2028
$ws.Range("N131").Value = -13629.4998
$ws.Range("H137").Value = 8101.1055
$ws.Range("I137").Value = 12972.8
$ws.Range("J137").Value = 2688.111
$ws.Range("K137").Value = 38918.39999999999
$ws.Range("L137").Value = 8064.333
$ws.Range("M137").Value = -33818.39999999999
$ws.Range("N137").Value = -18264.333
$ws.Range("H140").Value = 1514.7073
$ws.Range("J140").Value = 2865.3635
$ws.Range("L140").Value = 8596.0905
$ws.Range("N140").Value = -18956.0905

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 849.73334
$ws.Range("I107").Value = 544
$ws.Range("K107").Value = 544
$ws.Range("M107").Value = 1376
$ws.Range("H132").Value = 2664.8
$ws.Range("I132").Value = 1907.5
$ws.Range("J132").Value = 3530.2856
$ws.Range("K132").Value = 5722.5
$ws.Range("L132").Value = 10590.8568
$ws.Range("M132").Value = -3192.5
$ws.Range("N132").Value = -15650.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3644.8394
$ws.Range("I132").Value = 3651.6333
$ws.Range("J132").Value = 3637
$ws.Range("K132").Value = 10954.8999
$ws.Range("L132").Value = 10911
$ws.Range("M132").Value = -8424.8999
$ws.Range("N132").Value = -15971
